$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 833.3333
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("H4").Value = 439
$ws.Range("I4").Value = 439
$ws.Range("K4").Value = 439
$ws.Range("M4").Value = -325
$ws.Range("H17").Value = 820
$ws.Range("J17").Value = 820
$ws.Range("L17").Value = 2460
$ws.Range("N17").Value = -2796
$ws.Range("H40").Value = 2300
$ws.Range("I40").Value = 2250
$ws.Range("K40").Value = 2250
$ws.Range("M40").Value = -2075
$ws.Range("H51").Value = 6100
$ws.Range("J51").Value = 6100
$ws.Range("L51").Value = 6100
$ws.Range("N51").Value = -7068
$ws.Range("H76").Value = 3350635.5
$ws.Range("I76").Value = 4687371
$ws.Range("K76").Value = 4687371
$ws.Range("M76").Value = -4687056
$ws.Range("H79").Value = 3350635.5
$ws.Range("I79").Value = 4687371
$ws.Range("K79").Value = 4687371
$ws.Range("M79").Value = -4686279
$ws.Range("H98").Value = 2169
$ws.Range("I98").Value = 2344.35
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 2344.35
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = -846.3499999999999
$ws.Range("N98").Value = -3996
$ws.Range("H116").Value = 18228.428
$ws.Range("I116").Value = 35866.668
$ws.Range("J116").Value = 4999.75
$ws.Range("K116").Value = 35866.668
$ws.Range("L116").Value = 4999.75
$ws.Range("M116").Value = -32424.668
$ws.Range("N116").Value = -11883.75
$ws.Range("H121").Value = 787.5
$ws.Range("I121").Value = 575
$ws.Range("K121").Value = 1725
$ws.Range("M121").Value = 22
$ws.Range("H122").Value = 2169
$ws.Range("I122").Value = 2344.35
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 7033.049999999999
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -4583.049999999999
$ws.Range("N122").Value = -7900
$ws.Range("H137").Value = 31691.242
$ws.Range("I137").Value = 1268.32
$ws.Range("J137").Value = 126762.875
$ws.Range("K137").Value = 3804.96
$ws.Range("L137").Value = 380288.625
$ws.Range("M137").Value = -1254.96
$ws.Range("N137").Value = -385388.625
$ws.Range("H141").Value = 968012.4399999999
$ws.Range("I141").Value = 1219055.4
$ws.Range("K141").Value = 3657166.2
$ws.Range("M141").Value = -3651986.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 342.375
$ws.Range("I5").Value = 194.83333
$ws.Range("J5").Value = 785
$ws.Range("K5").Value = 194.83333
$ws.Range("L5").Value = 785
$ws.Range("M5").Value = -82.83332999999999
$ws.Range("N5").Value = -1009
$ws.Range("H61").Value = 5160.273
$ws.Range("I61").Value = 3416.5
$ws.Range("J61").Value = 7252.8
$ws.Range("K61").Value = 3416.5
$ws.Range("L61").Value = 7252.8
$ws.Range("M61").Value = -3204.5
$ws.Range("N61").Value = -7676.8
$ws.Range("H74").Value = 1092.8937
$ws.Range("I74").Value = 893.5526
$ws.Range("J74").Value = 1934.5555
$ws.Range("K74").Value = 893.5526
$ws.Range("L74").Value = 1934.5555
$ws.Range("M74").Value = -19.55259999999998
$ws.Range("N74").Value = -3682.5555
$ws.Range("H77").Value = 1092.8937
$ws.Range("I77").Value = 893.5526
$ws.Range("J77").Value = 1934.5555
$ws.Range("K77").Value = 4467.763
$ws.Range("L77").Value = 9672.7775
$ws.Range("M77").Value = -99.76299999999992
$ws.Range("N77").Value = -18408.7775
$ws.Range("H122").Value = 1616.75
$ws.Range("I122").Value = 1616.75
$ws.Range("K122").Value = 4850.25
$ws.Range("M122").Value = -2400.25
$ws.Range("H132").Value = 2352.318
$ws.Range("I132").Value = 2119.2222
$ws.Range("J132").Value = 2513.6924
$ws.Range("K132").Value = 6357.6666
$ws.Range("L132").Value = 7541.0772
$ws.Range("M132").Value = -3827.6666
$ws.Range("N132").Value = -12601.0772
$ws.Range("H136").Value = 5160.273
$ws.Range("I136").Value = 3416.5
$ws.Range("J136").Value = 7252.8
$ws.Range("K136").Value = 10249.5
$ws.Range("L136").Value = 21758.4
$ws.Range("M136").Value = -7699.5
$ws.Range("N136").Value = -26858.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 342.375
$ws.Range("I4").Value = 194.83333
$ws.Range("J4").Value = 785
$ws.Range("K4").Value = 194.83333
$ws.Range("L4").Value = 785
$ws.Range("M4").Value = -79.83332999999999
$ws.Range("N4").Value = -1015
$ws.Range("H105").Value = 2273.5454
$ws.Range("I105").Value = 2237.842
$ws.Range("K105").Value = 2237.842
$ws.Range("M105").Value = -490.8420000000001
$ws.Range("H134").Value = 7647.174
$ws.Range("I134").Value = 8836.579
$ws.Range("K134").Value = 26509.737
$ws.Range("M134").Value = -23974.737

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1864.3478
$ws.Range("I31").Value = 1103
$ws.Range("K31").Value = 1103
$ws.Range("M31").Value = -808
$ws.Range("H34").Value = 1864.3478
$ws.Range("I34").Value = 1103
$ws.Range("K34").Value = 1103
$ws.Range("M34").Value = -901

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 184.71428
$ws.Range("J2").Value = 149.5
$ws.Range("L2").Value = 897
$ws.Range("N2").Value = -1123
$ws.Range("H12").Value = 250.25
$ws.Range("J12").Value = 250.25
$ws.Range("L12").Value = 750.75
$ws.Range("N12").Value = -1096.75
$ws.Range("H38").Value = 381.5
$ws.Range("I38").Value = 65
$ws.Range("J38").Value = 698
$ws.Range("K38").Value = 195
$ws.Range("L38").Value = 2094
$ws.Range("M38").Value = 152
$ws.Range("N38").Value = -2788
$ws.Range("H97").Value = 963.25
$ws.Range("I97").Value = 940
$ws.Range("K97").Value = 2820
$ws.Range("M97").Value = -2324
$ws.Range("H122").Value = 912.38464
$ws.Range("J122").Value = 1155.75
$ws.Range("L122").Value = 10401.75
$ws.Range("N122").Value = -15301.75
$ws.Range("H131").Value = 788.4299999999999
$ws.Range("J131").Value = 799.40625
$ws.Range("L131").Value = 2398.21875
$ws.Range("N131").Value = -12478.21875
$ws.Range("H139").Value = 21800
$ws.Range("J139").Value = 2250
$ws.Range("L139").Value = 6750
$ws.Range("N139").Value = -17030
$ws.Range("H140").Value = 1635.0278
$ws.Range("I140").Value = 932.41174
$ws.Range("K140").Value = 2797.23522
$ws.Range("M140").Value = 2382.76478
$ws.Range("H141").Value = 2993.2856
$ws.Range("I141").Value = 2993.2856
$ws.Range("K141").Value = 8979.856800000001
$ws.Range("M141").Value = -3799.856800000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2486.24
$ws.Range("I102").Value = 2557.889
$ws.Range("J102").Value = 2302
$ws.Range("K102").Value = 2557.889
$ws.Range("L102").Value = 2302
$ws.Range("M102").Value = -935.8890000000001
$ws.Range("N102").Value = -5546
$ws.Range("H122").Value = 1711.5
$ws.Range("I122").Value = 1291.8182
$ws.Range("K122").Value = 3875.4546
$ws.Range("M122").Value = -1425.4546
$ws.Range("H132").Value = 3208036.5
$ws.Range("I132").Value = 7694252
$ws.Range("K132").Value = 23082756
$ws.Range("M132").Value = -23080226

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1696.9474
$ws.Range("I132").Value = 1404.1765
$ws.Range("J132").Value = 1933.9524
$ws.Range("K132").Value = 4212.529500000001
$ws.Range("L132").Value = 5801.857199999999
$ws.Range("M132").Value = -1682.529500000001
$ws.Range("N132").Value = -10861.8572

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 66583.25
$ws.Range("I122").Value = 156775
$ws.Range("J122").Value = 2160.5715
$ws.Range("K122").Value = 470325
$ws.Range("L122").Value = 6481.7145
$ws.Range("M122").Value = -467875
$ws.Range("N122").Value = -11381.7145
